# Cleaned up data and added more detailed benchmarking to the dashboard.
# Rewrites the Benchmark_Data_Pivot sheet data (A2:C101): each of the 5 categories
# now has 20 rows (Index 1-20) instead of 14, with refreshed Score values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  @(1, "Work Life Balance", 67.20473385133047),
  @(2, "Work Life Balance", 78.63010069662232),
  @(3, "Work Life Balance", 72.43458032973052),
  @(4, "Work Life Balance", 29.81239756343025),
  @(5, "Work Life Balance", 112.467889172675),
  @(6, "Work Life Balance", 47.97833772166527),
  @(7, "Work Life Balance", 115.1250582621385),
  @(8, "Work Life Balance", 80.79870948563672),
  @(9, "Work Life Balance", 74.95113270080171),
  @(10, "Work Life Balance", 60.77557172992015),
  @(11, "Work Life Balance", 97.76304832841411),
  @(12, "Work Life Balance", 109.235535057134),
  @(13, "Work Life Balance", 43.74138976081582),
  @(14, "Work Life Balance", 37.5766571147183),
  @(15, "Work Life Balance", 30.19232658888301),
  @(16, "Work Life Balance", 35.40293979630021),
  @(17, "Work Life Balance", 108.0975997812136),
  @(18, "Work Life Balance", 96.56955052210475),
  @(19, "Work Life Balance", 28.54760002264616),
  @(20, "Work Life Balance", 100.2466061084747),
  @(1, "Career Growth", 46.07665320990305),
  @(2, "Career Growth", 73.33352598562178),
  @(3, "Career Growth", 36.39006831610373),
  @(4, "Career Growth", 21.56964384126162),
  @(5, "Career Growth", 72.9822244490293),
  @(6, "Career Growth", 94.406976086239),
  @(7, "Career Growth", 100.0408759846802),
  @(8, "Career Growth", 20.62052569120796),
  @(9, "Career Growth", 33.76300632014023),
  @(10, "Career Growth", 99.69559907193678),
  @(11, "Career Growth", 56.57585366599079),
  @(12, "Career Growth", 47.13436980365095),
  @(13, "Career Growth", 79.53369367425483),
  @(14, "Career Growth", 34.64945386748243),
  @(15, "Career Growth", 23.70128457989871),
  @(16, "Career Growth", 50.47987470739064),
  @(17, "Career Growth", 32.53923203854022),
  @(18, "Career Growth", 106.9302808920604),
  @(19, "Career Growth", 61.7410190091136),
  @(20, "Career Growth", 40.82497456475512),
  @(1, "Compensation", 52.85596913433723),
  @(2, "Compensation", 41.33022577225088),
  @(3, "Compensation", 43.51467424933426),
  @(4, "Compensation", 53.42496593739565),
  @(5, "Compensation", 37.93595952975845),
  @(6, "Compensation", 123.8525871669547),
  @(7, "Compensation", 124.1310424330165),
  @(8, "Compensation", 96.23963159904581),
  @(9, "Compensation", 65.70514376943694),
  @(10, "Compensation", 75.81106431555679),
  @(11, "Compensation", 25.36446735226679),
  @(12, "Compensation", 106.303752071009),
  @(13, "Compensation", 38.77805462161064),
  @(14, "Compensation", 32.84387486945073),
  @(15, "Compensation", 72.7759513777151),
  @(16, "Compensation", 92.04112859555397),
  @(17, "Compensation", 53.24497978772666),
  @(18, "Compensation", 97.81915376489032),
  @(19, "Compensation", 36.23125820436501),
  @(20, "Compensation", 67.72469838828945),
  @(1, "Leadership", 42.60660710848016),
  @(2, "Leadership", 90.97811146643561),
  @(3, "Leadership", 76.81837814287593),
  @(4, "Leadership", 57.28534498478044),
  @(5, "Leadership", 43.0134711060215),
  @(6, "Leadership", 53.43798443339365),
  @(7, "Leadership", 16.32924520076737),
  @(8, "Leadership", 102.8259658562364),
  @(9, "Leadership", 45.41043831657473),
  @(10, "Leadership", 70.98309835611551),
  @(11, "Leadership", 18.96447130779364),
  @(12, "Leadership", 64.65109897466549),
  @(13, "Leadership", 50.67380076399056),
  @(14, "Leadership", 110.9899305992184),
  @(15, "Leadership", 39.03099327094822),
  @(16, "Leadership", 56.7396242537396),
  @(17, "Leadership", 80.83670967033808),
  @(18, "Leadership", 65.79838219309795),
  @(19, "Leadership", 104.2436916048489),
  @(20, "Leadership", 34.6735428564058),
  @(1, "Colleagues", 65.66261747589319),
  @(2, "Colleagues", 37.94429338693888),
  @(3, "Colleagues", 75.45115842914211),
  @(4, "Colleagues", 75.75407140709376),
  @(5, "Colleagues", 65.08766548335271),
  @(6, "Colleagues", 109.4509584341707),
  @(7, "Colleagues", 76.54758765971221),
  @(8, "Colleagues", 62.56582209239682),
  @(9, "Colleagues", 109.0581330574912),
  @(10, "Colleagues", 50.69110425082442),
  @(11, "Colleagues", 35.09838467933363),
  @(12, "Colleagues", 127.6185571242596),
  @(13, "Colleagues", 72.57105192957494),
  @(14, "Colleagues", 84.06646522653472),
  @(15, "Colleagues", 80.6410620130309),
  @(16, "Colleagues", 72.769955636124),
  @(17, "Colleagues", 74.17440152235834),
  @(18, "Colleagues", 67.95297459263324),
  @(19, "Colleagues", 90.15648100482372),
  @(20, "Colleagues", 70.06641566555348)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
